$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$school = "โรงเรียนโคกเพชรวิทยาคาร"
$detailEdu = "วัสดุการศึกษา"
$detailOffice = "วัสดุสำนักงาน"
$detailEduCost = "ค่าวัสดุการศึกษา"
$dateStr = "Created on 14-01-2024"

$rows = @(
    @(25, $school, $detailEdu,   1, $dateStr),
    @(26, $school, $detailEdu,   1, $dateStr),
    @(27, $school, $detailEdu,   1, $dateStr),
    @(28, $school, $detailEdu,   1, $dateStr),
    @(29, $school, $detailEdu,   1, $dateStr),
    @(30, $school, $detailOffice, 2, $dateStr),
    @(31, $school, $detailEdu,   1, $dateStr),
    @(32, $school, $detailEdu,   1, $dateStr),
    @(33, $school, $detailEdu,   1, $dateStr),
    @(34, $school, $detailEdu,   1, $dateStr),
    @(35, $school, $detailEduCost, 5, $dateStr)
)

$r = 29
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
